$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, copy the current (soon to be old) row 14 values down to the new row 15,
# since this existing weekly entry is being kept and a new entry is being
# recorded in row 14 for the newer date.
$ws.Range("A15").Value = $ws.Range("A14").Value2
$ws.Range("B15").Value = $ws.Range("B14").Value2
$ws.Range("C15").Value = $ws.Range("C14").Value2
$ws.Range("D15").Value = $ws.Range("D14").Value2
$ws.Range("D15").NumberFormat = $ws.Range("D14").NumberFormat
$ws.Range("E15").Value = $ws.Range("E14").Value2
$ws.Range("F15").Value = $ws.Range("F14").Value2
$ws.Range("G15").Value = $ws.Range("G14").Value2
$ws.Range("H15").Value = $ws.Range("H14").Value2
$ws.Range("I15").Value = $ws.Range("I14").Value2
$ws.Range("J15").Value = $ws.Range("J14").Value2
$ws.Range("K15").Value = $ws.Range("K14").Value2
$ws.Range("L15").Value = $ws.Range("L14").Value2
$ws.Range("M15").Value = $ws.Range("M14").Value2
$ws.Range("N15").Value = $ws.Range("N14").Value2
$ws.Range("O15").Value = $ws.Range("O14").Value2
$ws.Range("P15").Value = $ws.Range("P14").Value2
$ws.Range("Q15").Value = $ws.Range("Q14").Value2
$ws.Range("R15").Value = $ws.Range("R14").Value2

# Now update row 14 with the new weekly values.
$ws.Range("D14").Value = 44474
$ws.Range("J14").Value = 18
$ws.Range("K14").Value = 100000
$ws.Range("L14").Value = 100000
$ws.Range("M14").Value = 100000
$ws.Range("P14").Value = 4000
